$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D, shifting existing D:K data to E:L
$ws.Columns.Item(4).Insert()

# Copy number formatting/style from the (now shifted) column E into the new column D
# (limited to the used range only, to avoid materializing the entire 1M+ row column)
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Populate the new column D with the latest period's financial figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 102100
$ws.Range("D15").Value = -600
$ws.Range("D17").Value = 26700
$ws.Range("D18").Value = 75400
$ws.Range("D20").Value = -29200
$ws.Range("D21").Value = 47900
$ws.Range("D23").Value = 46200
$ws.Range("D24").Value = 10600
$ws.Range("D26").Value = 35600
$ws.Range("D27").Value = 35600
$ws.Range("D29").Value = 500
$ws.Range("D32").Value = 29200
$ws.Range("D33").Value = 36100
$ws.Range("D35").Value = 36100
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 147700
$ws.Range("D42").Value = 10300
$ws.Range("D48").Value = 17300
$ws.Range("D49").Value = 66000
$ws.Range("D52").Value = 4600
$ws.Range("D54").Value = 2974000
$ws.Range("D57").Value = 21900
$ws.Range("D59").Value = "NA"
$ws.Range("D61").Value = 113200
$ws.Range("D66").Value = 2599500
$ws.Range("D72").Value = 81600
$ws.Range("D76").Value = 374500
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 36100
$ws.Range("D83").Value = 1700
$ws.Range("D89").Value = -84600
$ws.Range("D91").Value = -2500
$ws.Range("D94").Value = -322400
$ws.Range("D100").Value = 404700
$ws.Range("D102").Value = -2400

# Rows where every period column already reads "NA" - the new column keeps the same "NA" marker
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"

# Rows where every period column already reads 0 - the new column keeps the same 0 value
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D77").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D101").Value = 0
